$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1011.7143
$ws.Cells.Item(19, 10).Value = 1013.6667
$ws.Cells.Item(19, 12).Value = 1013.6667
$ws.Cells.Item(19, 14).Value = -1363.6667

$ws.Cells.Item(33, 8).Value = 217.5625
$ws.Cells.Item(33, 9).Value = 126.75
$ws.Cells.Item(33, 10).Value = 490
$ws.Cells.Item(33, 11).Value = 126.75
$ws.Cells.Item(33, 12).Value = 490
$ws.Cells.Item(33, 13).Value = 102.25
$ws.Cells.Item(33, 14).Value = -948

$ws.Cells.Item(38, 8).Value = 362.125
$ws.Cells.Item(38, 9).Value = 149.5
$ws.Cells.Item(38, 11).Value = 448.5
$ws.Cells.Item(38, 13).Value = -76.5

$ws.Cells.Item(69, 8).Value = 3015
$ws.Cells.Item(69, 10).Value = 3015
$ws.Cells.Item(69, 12).Value = 9045
$ws.Cells.Item(69, 14).Value = -10793

$ws.Cells.Item(72, 8).Value = 3015
$ws.Cells.Item(72, 10).Value = 3015
$ws.Cells.Item(72, 12).Value = 27135
$ws.Cells.Item(72, 14).Value = -35871

$ws.Cells.Item(76, 8).Value = 8187.4
$ws.Cells.Item(76, 9).Value = 12190.909
$ws.Cells.Item(76, 10).Value = 3294.2222
$ws.Cells.Item(76, 11).Value = 12190.909
$ws.Cells.Item(76, 12).Value = 3294.2222
$ws.Cells.Item(76, 13).Value = -11875.909
$ws.Cells.Item(76, 14).Value = -3924.2222

$ws.Cells.Item(79, 8).Value = 8187.4
$ws.Cells.Item(79, 9).Value = 12190.909
$ws.Cells.Item(79, 10).Value = 3294.2222
$ws.Cells.Item(79, 11).Value = 12190.909
$ws.Cells.Item(79, 12).Value = 3294.2222
$ws.Cells.Item(79, 13).Value = -11098.909
$ws.Cells.Item(79, 14).Value = -5478.2222

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 13).ClearContents()

$ws.Cells.Item(32, 8).Value = 21532.03
$ws.Cells.Item(32, 9).Value = 19429.684
$ws.Cells.Item(32, 11).Value = 19429.684
$ws.Cells.Item(32, 13).Value = -19142.684

$ws.Cells.Item(37, 8).Value = 35722.223
$ws.Cells.Item(37, 10).Value = 40000
$ws.Cells.Item(37, 12).Value = 40000
$ws.Cells.Item(37, 14).Value = -40546

$ws.Cells.Item(55, 8).Value = 11500
$ws.Cells.Item(55, 9).Value = 3000
$ws.Cells.Item(55, 11).Value = 3000
$ws.Cells.Item(55, 13).Value = -2685

$ws.Cells.Item(80, 8).Value = 34935
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 34935
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 34935
$ws.Cells.Item(80, 14).Value = -36931
$ws.Cells.Item(80, 13).ClearContents()

$ws.Cells.Item(83, 8).Value = 34935
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 34935
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 104805
$ws.Cells.Item(83, 14).Value = -114789
$ws.Cells.Item(83, 13).ClearContents()

$ws.Cells.Item(132, 8).Value = 42786.715
$ws.Cells.Item(132, 9).Value = 26523.9
$ws.Cells.Item(132, 10).Value = 115065.89
$ws.Cells.Item(132, 11).Value = 79571.70000000001
$ws.Cells.Item(132, 12).Value = 345197.67
$ws.Cells.Item(132, 13).Value = -77041.70000000001
$ws.Cells.Item(132, 14).Value = -350257.67

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 321.35715
$ws.Cells.Item(22, 9).Value = 208.33333
$ws.Cells.Item(22, 11).Value = 208.33333
$ws.Cells.Item(22, 13).Value = -35.33332999999999

$ws.Cells.Item(105, 8).Value = 50302220
$ws.Cells.Item(105, 9).Value = 62877036
$ws.Cells.Item(105, 11).Value = 62877036
$ws.Cells.Item(105, 13).Value = -62875289

$ws.Cells.Item(134, 8).Value = 2889
$ws.Cells.Item(134, 9).Value = 3025.125
$ws.Cells.Item(134, 11).Value = 9075.375
$ws.Cells.Item(134, 13).Value = -6540.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1476
$ws.Cells.Item(16, 9).Value = 920.8333
$ws.Cells.Item(16, 10).Value = 2142.2
$ws.Cells.Item(16, 11).Value = 920.8333
$ws.Cells.Item(16, 12).Value = 2142.2
$ws.Cells.Item(16, 13).Value = -633.8333
$ws.Cells.Item(16, 14).Value = -2716.2

$ws.Cells.Item(113, 8).Value = 1476
$ws.Cells.Item(113, 9).Value = 920.8333
$ws.Cells.Item(113, 10).Value = 2142.2
$ws.Cells.Item(113, 11).Value = 920.8333
$ws.Cells.Item(113, 12).Value = 2142.2
$ws.Cells.Item(113, 13).Value = 1249.1667
$ws.Cells.Item(113, 14).Value = -6482.2

$ws.Cells.Item(132, 8).Value = 52648.15
$ws.Cells.Item(132, 9).Value = 2668.4614
$ws.Cells.Item(132, 10).Value = 145467.58
$ws.Cells.Item(132, 11).Value = 8005.3842
$ws.Cells.Item(132, 12).Value = 436402.74
$ws.Cells.Item(132, 13).Value = -5475.3842
$ws.Cells.Item(132, 14).Value = -441462.74

$ws.Cells.Item(134, 8).Value = 62643.277
$ws.Cells.Item(134, 9).Value = 2954
$ws.Cells.Item(134, 10).Value = 217835.4
$ws.Cells.Item(134, 11).Value = 8862
$ws.Cells.Item(134, 12).Value = 653506.2
$ws.Cells.Item(134, 13).Value = -6327
$ws.Cells.Item(134, 14).Value = -658576.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(43, 8).Value = 4117
$ws.Cells.Item(43, 10).Value = 4380
$ws.Cells.Item(43, 12).Value = 13140
$ws.Cells.Item(43, 14).Value = -13368

$ws.Cells.Item(131, 8).Value = 1200.6809
$ws.Cells.Item(131, 9).Value = 626.25
$ws.Cells.Item(131, 10).Value = 1497.1613
$ws.Cells.Item(131, 11).Value = 1878.75
$ws.Cells.Item(131, 12).Value = 4491.4839
$ws.Cells.Item(131, 13).Value = 3161.25
$ws.Cells.Item(131, 14).Value = -14571.4839

$ws.Cells.Item(132, 8).Value = 1386.375
$ws.Cells.Item(132, 10).Value = 1962.875
$ws.Cells.Item(132, 12).Value = 17665.875
$ws.Cells.Item(132, 14).Value = -22725.875

$ws.Cells.Item(141, 8).Value = 6440.5557
$ws.Cells.Item(141, 9).Value = 3129.2307
$ws.Cells.Item(141, 10).Value = 15050
$ws.Cells.Item(141, 11).Value = 9387.6921
$ws.Cells.Item(141, 12).Value = 45150
$ws.Cells.Item(141, 13).Value = -4207.6921
$ws.Cells.Item(141, 14).Value = -55510

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 171.33333
$ws.Cells.Item(2, 9).Value = 134
$ws.Cells.Item(2, 10).Value = 246
$ws.Cells.Item(2, 11).Value = 134
$ws.Cells.Item(2, 12).Value = 246
$ws.Cells.Item(2, 13).Value = -21
$ws.Cells.Item(2, 14).Value = -472

$ws.Cells.Item(70, 8).Value = 32110.621
$ws.Cells.Item(70, 9).Value = 50313.91
$ws.Cells.Item(70, 10).Value = 5412.467
$ws.Cells.Item(70, 11).Value = 50313.91
$ws.Cells.Item(70, 12).Value = 5412.467
$ws.Cells.Item(70, 13).Value = -50043.91
$ws.Cells.Item(70, 14).Value = -5952.467

$ws.Cells.Item(73, 8).Value = 32110.621
$ws.Cells.Item(73, 9).Value = 50313.91
$ws.Cells.Item(73, 10).Value = 5412.467
$ws.Cells.Item(73, 11).Value = 50313.91
$ws.Cells.Item(73, 12).Value = 5412.467
$ws.Cells.Item(73, 13).Value = -49377.91
$ws.Cells.Item(73, 14).Value = -7284.467

$ws.Cells.Item(132, 8).Value = 81876.88
$ws.Cells.Item(132, 9).Value = 84235.5
$ws.Cells.Item(132, 11).Value = 252706.5
$ws.Cells.Item(132, 13).Value = -250176.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 1918
$ws.Cells.Item(68, 9).Value = 1493.3334
$ws.Cells.Item(68, 10).Value = 2100
$ws.Cells.Item(68, 11).Value = 1493.3334
$ws.Cells.Item(68, 12).Value = 2100
$ws.Cells.Item(68, 13).Value = -744.3334
$ws.Cells.Item(68, 14).Value = -3598

$ws.Cells.Item(71, 8).Value = 1918
$ws.Cells.Item(71, 9).Value = 1493.3334
$ws.Cells.Item(71, 10).Value = 2100
$ws.Cells.Item(71, 11).Value = 7466.666999999999
$ws.Cells.Item(71, 12).Value = 10500
$ws.Cells.Item(71, 13).Value = -3722.666999999999
$ws.Cells.Item(71, 14).Value = -17988

$ws.Cells.Item(93, 8).Value = 300
$ws.Cells.Item(93, 9).Value = 300
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 300
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = 948
$ws.Cells.Item(93, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 52748.15
$ws.Cells.Item(132, 9).Value = 1138.6666
$ws.Cells.Item(132, 10).Value = 130162.375
$ws.Cells.Item(132, 11).Value = 3415.9998
$ws.Cells.Item(132, 12).Value = 390487.125
$ws.Cells.Item(132, 13).Value = -885.9998000000001
$ws.Cells.Item(132, 14).Value = -395547.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4445.4165
$ws.Cells.Item(62, 9).Value = 4000
$ws.Cells.Item(62, 10).Value = 5069
$ws.Cells.Item(62, 11).Value = 4000
$ws.Cells.Item(62, 12).Value = 5069
$ws.Cells.Item(62, 13).Value = -3376
$ws.Cells.Item(62, 14).Value = -6317

$ws.Cells.Item(65, 8).Value = 4445.4165
$ws.Cells.Item(65, 9).Value = 4000
$ws.Cells.Item(65, 10).Value = 5069
$ws.Cells.Item(65, 11).Value = 20000
$ws.Cells.Item(65, 12).Value = 25345
$ws.Cells.Item(65, 13).Value = -16880
$ws.Cells.Item(65, 14).Value = -31585

$ws.Cells.Item(132, 8).Value = 60234.176
$ws.Cells.Item(132, 9).Value = 42679.332
$ws.Cells.Item(132, 10).Value = 102365.8
$ws.Cells.Item(132, 11).Value = 128037.996
$ws.Cells.Item(132, 12).Value = 307097.4
$ws.Cells.Item(132, 13).Value = -125507.996
$ws.Cells.Item(132, 14).Value = -312157.4
